# Weekly update: insert a new Albahaca price record at row 11 of Sheet1,
# pushing the existing records (old rows 11-131) down by one row
# (new rows 12-132). Dimension grows from A1:R131 to A1:R132.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 11; Excel shifts rows 11..131 down to
# 12..132 and copies the row-above formatting (e.g. the date number format
# on column D) into the newly inserted, still-empty row.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new weekly record.
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44750
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = 100112052
$ws.Range("G11").Value = "Albahaca"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 6000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 6500
$ws.Range("N11").Value = "`$/paquete"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 6500
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
